$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- W1: new header cell "08-10-2020", styled like V1 (bold/center/top, thin border) ---
# The leading apostrophe forces a text entry so Excel does not auto-convert the
# date-like string into a date serial number (matches the inlineStr text in the diff).
$srcHeader = $ws.Range("V1")
$dstHeader = $ws.Range("W1")
$dstHeader.Value = "'08-10-2020"
$dstHeader.Font.Name = $srcHeader.Font.Name
$dstHeader.Font.Bold = $srcHeader.Font.Bold
$dstHeader.Font.Size = $srcHeader.Font.Size
$dstHeader.Borders.LineStyle = $srcHeader.Borders.Item(7).LineStyle
$dstHeader.HorizontalAlignment = $srcHeader.HorizontalAlignment
$dstHeader.VerticalAlignment = $srcHeader.VerticalAlignment

# --- W2:W36: new daily confirmed-case totals for 08-10-2020 ---
$ws.Range("W2").Value = 3696
$ws.Range("W3").Value = 678828
$ws.Range("W4").Value = 8396
$ws.Range("W5").Value = 157638
$ws.Range("W6").Value = 179732
$ws.Range("W7").Value = 11190
$ws.Range("W8").Value = 103828
$ws.Range("W9").Value = 3010
$ws.Range("W10").Value = 270305
$ws.Range("W11").Value = 31444
$ws.Range("W12").Value = 126657
$ws.Range("W13").Value = 124841
$ws.Range("W14").Value = 13338
$ws.Range("W15").Value = 67684
$ws.Range("W16").Value = 79176
$ws.Range("W17").Value = 542906
$ws.Range("W18").Value = 160253
$ws.Range("W19").Value = 3511
$ws.Range("W20").Value = 120267
$ws.Range("W21").Value = 1196441
$ws.Range("W22").Value = 9604
$ws.Range("W23").Value = 4694
$ws.Range("W24").Value = 1919
$ws.Range("W25").Value = 5498
$ws.Range("W26").Value = 213672
$ws.Range("W27").Value = 24930
$ws.Range("W28").Value = 105585
$ws.Range("W29").Value = 127526
$ws.Range("W30").Value = 2615
$ws.Range("W31").Value = 580736
$ws.Range("W32").Value = 179075
$ws.Range("W33").Value = 23066
$ws.Range("W34").Value = 43904
$ws.Range("W35").Value = 374972
$ws.Range("W36").Value = 246767

